$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6792637705802917
$ws.Range("B1").Value = 3.302244901657104
$ws.Range("C1").Value = 2.791783571243286
$ws.Range("D1").Value = 1.053259134292603
$ws.Range("E1").Value = 1.016553401947021
